$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New daily collection entries (column P = new date column) and a few
# --- other in-place corrections for the "till 12 Dec 9pm" update ---

$ws.Range("P4").Value = 3000
$ws.Range("P7").Value = 1000
$ws.Range("P21").Value = 1000
$ws.Range("P28").Value = 100
$ws.Range("P39").Value = 3000
$ws.Range("P53").Value = 5000
$ws.Range("K54").Value = 3000
$ws.Range("P54").Value = 3000
$ws.Range("P61").Value = 1000
$ws.Range("I62").Value = 2000
$ws.Range("B62").Value = "WALIDAD"
$ws.Range("P65").Value = 1000
$ws.Range("P67").Value = 500
$ws.Range("P71").Value = 1000
$ws.Range("P83").Value = 3000
$ws.Range("P84").Value = 3000

# Row 24's collection entry was moved/removed - clear the old J24 figure
$ws.Range("J24").Clear()

# Tidy up: drop the stray blank "Area" placeholder cells in column B that
# don't carry a value (so the cell itself no longer exists in the sheet)
$blankAreaRows = 11,12,15,17,19,20,22,24,26,31,32,37,41,42,44,45,48,49,52,55,56,57,58,60,69,72,73,75,76,81
foreach ($r in $blankAreaRows) {
    $ws.Range("B$r").Clear()
}

# Refresh the frozen panes so column D is included with the header rows,
# and leave the view focused on the most recently edited area
$ws.Range("A1").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("E3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F8").Select()
